$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format so numeric-looking price strings
# (e.g. "206.70", "1.00", "0.690") are preserved exactly as text
# instead of being auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '87.476.95'
$ws.Range("E2").Value = '  +0.23%  '

# Row 3
$ws.Range("D3").Value = '3.158.75'
$ws.Range("E3").Value = '  -5.49%  '

# Row 4
$ws.Range("E4").Value = '  -0.33%  '

# Row 5
$ws.Range("D5").Value = '206.70'
$ws.Range("E5").Value = '  -3.51%  '

# Row 6
$ws.Range("D6").Value = '606.96'
$ws.Range("E6").Value = '  -3.49%  '

# Row 7
$ws.Range("D7").Value = '0.378'
$ws.Range("E7").Value = '  -4.36%  '

# Row 8
$ws.Range("D8").Value = '0.662'
$ws.Range("E8").Value = '  +3.75%  '

# Row 9
$ws.Range("D9").Value = '1.00'
$ws.Range("E9").Value = '  -0.22%  '

# Row 10
$ws.Range("D10").Value = '3.159.82'
$ws.Range("E10").Value = '  -5.88%  '

# Row 11
$ws.Range("D11").Value = '0.532'
$ws.Range("E11").Value = '  -10.57%  '

# Row 12
$ws.Range("E12").Value = '  +5.92%  '

# Row 13
$ws.Range("E13").Value = '  -10.07%  '

# Row 14
$ws.Range("D14").Value = '3.743.94'
$ws.Range("E14").Value = '  -5.92%  '

# Row 15
$ws.Range("D15").Value = '5.21'
$ws.Range("E15").Value = '  -2.63%  '

# Row 16
$ws.Range("D16").Value = '87.275.62'
$ws.Range("E16").Value = '  -0.21%  '

# Row 17
$ws.Range("D17").Value = '31.89'
$ws.Range("E17").Value = '  -9.13%  '

# Row 18
$ws.Range("D18").Value = '3.160.86'
$ws.Range("E18").Value = '  -6.30%  '

# Row 19
$ws.Range("D19").Value = '3.07'
$ws.Range("E19").Value = '  +4.80%  '

# Row 20
$ws.Range("D20").Value = '13.30'
$ws.Range("E20").Value = '  -8.27%  '

# Row 21
$ws.Range("D21").Value = '411.73'
$ws.Range("E21").Value = '  -6.91%  '

# Row 22
$ws.Range("D22").Value = '8.38'
$ws.Range("E22").Value = '  -9.89%  '

# Row 23
$ws.Range("D23").Value = '5.01'
$ws.Range("E23").Value = '  -7.12%  '

# Row 24
$ws.Range("D24").Value = '5.12'
$ws.Range("E24").Value = '  -3.29%  '

# Row 25
$ws.Range("D25").Value = '11.78'
$ws.Range("E25").Value = '  -2.91%  '

# Row 26
$ws.Range("D26").Value = '3.331.39'
$ws.Range("E26").Value = '  -6.44%  '

# Row 27
$ws.Range("D27").Value = '72.88'
$ws.Range("E27").Value = '  -7.40%  '

# Row 28
$ws.Range("D28").Value = '0.0000129'
$ws.Range("E28").Value = '  -3.97%  '

# Row 29
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  +0.02%  '

# Row 30
$ws.Range("B30").Value = 'Cronos'
$ws.Range("C30").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D30").Value = '0.158'
$ws.Range("E30").Value = '  -9.31%  '

# Row 31
$ws.Range("B31").Value = 'Binance-PegBSC-USD'
$ws.Range("C31").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D31").Value = '0.999'
$ws.Range("E31").Value = '  -0.22%  '

# Row 32
$ws.Range("D32").Value = '539.48'
$ws.Range("E32").Value = '  -2.77%  '

# Row 33
$ws.Range("D33").Value = '8.16'
$ws.Range("E33").Value = '  -9.60%  '

# Row 34
$ws.Range("E34").Value = '  -13.40%  '

# Row 35
$ws.Range("D35").Value = '1.83'
$ws.Range("E35").Value = '  -9.24%  '

# Row 36
$ws.Range("D36").Value = '6.59'
$ws.Range("E36").Value = '  -5.25%  '

# Row 37
$ws.Range("D37").Value = '0.131'
$ws.Range("E37").Value = '  -4.49%  '

# Row 38
$ws.Range("D38").Value = '21.68'
$ws.Range("E38").Value = '  -5.97%  '

# Row 39
$ws.Range("D39").Value = '21.82'
$ws.Range("E39").Value = '  -0.02%  '

# Row 40
$ws.Range("E40").Value = '  -0.16%  '

# Row 41
$ws.Range("D41").Value = '2.98'
$ws.Range("E41").Value = '  +0.65%  '

# Row 42
$ws.Range("E42").Value = '  +0.15%  '

# Row 43
$ws.Range("E43").Value = '  -6.74%  '

# Row 44
$ws.Range("D44").Value = '0.364'
$ws.Range("E44").Value = '  -12.67%  '

# Row 45
$ws.Range("D45").Value = '148.91'
$ws.Range("E45").Value = '  -5.83%  '

# Row 46
$ws.Range("D46").Value = '171.52'
$ws.Range("E46").Value = '  -5.87%  '

# Row 47
$ws.Range("D47").Value = '43.06'
$ws.Range("E47").Value = '  -6.52%  '

# Row 48
$ws.Range("E48").Value = '  +2.36%  '

# Row 49
$ws.Range("D49").Value = '1.21'
$ws.Range("E49").Value = '  -11.43%  '

# Row 50
$ws.Range("D50").Value = '3.93'
$ws.Range("E50").Value = '  -9.83%  '

# Row 51
$ws.Range("B51").Value = 'Mantle'
$ws.Range("C51").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D51").Value = '0.690'
$ws.Range("E51").Value = '  -10.47%  '
